$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (LDA / exp)
$ws.Range("C2").Value = 0.008981458333333329
$ws.Range("D2").Value = 118.5795307831199
$ws.Range("E2").Value = 0.01079457470020627
$ws.Range("F2").Value = 0.008981458333333329

# Row 3 (LDA / ethanol)
$ws.Range("C3").Value = 0.008981458333333331
$ws.Range("D3").Value = 118.5795307831199
$ws.Range("E3").Value = 0.01079457470020627
$ws.Range("F3").Value = 0.008981458333333331

# Row 4 (PBE / ethanol)
$ws.Range("C4").Value = 0.008981458333333329
$ws.Range("D4").Value = 118.5795307831199
$ws.Range("E4").Value = 0.01079457470020627
$ws.Range("F4").Value = 0.008981458333333329
